$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the full row contents (columns B:K) between row 4 (Deandre Ayton) and
# row 5 (Jock Landale), while leaving column A ("No.") untouched. This moves
# Jock Landale up to roster position 2 and Deandre Ayton down to position 3.

$row4 = $ws.Range("B4:K4").Value()
$row5 = $ws.Range("B5:K5").Value()

$ws.Range("B4:K4").Value = $row5
$ws.Range("B5:K5").Value = $row4
